$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(93).Insert()

$ws.Range("A93").Value = 5
$ws.Range("B93").Value = "Macroferia Regional de Talca"
$ws.Range("C93").Value = "Maule"
$ws.Range("D93").Value = 44741
$ws.Range("E93").Value = 7
$ws.Range("F93").Value = 100112031
$ws.Range("G93").Value = "Poroto verde"
$ws.Range("H93").Value = "Sin especificar"
$ws.Range("I93").Value = "Primera"
$ws.Range("J93").Value = 150
$ws.Range("K93").Value = 25000
$ws.Range("L93").Value = 25000
$ws.Range("M93").Value = 25000
$ws.Range("N93").Value = "$/malla 25 kilos"
$ws.Range("O93").Value = "Región de Arica y Parinacota"
$ws.Range("P93").Value = 1000
$ws.Range("Q93").Value = 25
$ws.Range("R93").Value = "Hortaliza"
